# Revisi Report (Bissmillah bnr), add Notification Import, Revisi Import
#
# - Rename header "EmployeeID" (B1) -> "Employee ID"
# - Rename header "Balance" (K1) -> "Amount"
# - Highlight the import header row (A1:K1) with a green/white "notification"
#   style: bold white text, forest-green fill, thin black border, centered.
# - Move the active selection to the header row A1:K1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text revisions -------------------------------------------------
$ws.Range("B1").Value = "Employee ID"
$ws.Range("K1").Value = "Amount"

# --- New "import notification" header style for A1:K1 ----------------------
$rng = $ws.Range("A1:K1")

# Font: bold, white, explicit Calibri (no theme scheme)
$rng.Font.Name = "Calibri"
$rng.Font.Bold = $true
$rng.Font.Color = 16777215

# Border: thin black box around every header cell
$rng.Borders.Color = 0
$rng.Borders.LineStyle = 1

# Fill: forest green interior
$rng.Interior.PatternColor = 0
$rng.Interior.Color = 2263842

# Alignment: centered both ways
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# --- Update the selection to match the highlighted header range ------------
$ws.Range("A1:K1").Select()
